# Add the new data row (case 6308) to the Optical_Power tracking sheet,
# matching the existing rows' layout: columns A-L are plain text
# (even the digit-only / date-looking values), M & N are numeric
# coordinates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1
if ($row -lt 2) { $row = 2 }

# Force columns A:L to be stored as text so values like "6308", "6", "1"
# and the date string "7/1/2025" are not reinterpreted as numbers/dates
# (mirrors how the existing rows were authored as inline strings).
$textRange = $ws.Range("A" + $row + ":L" + $row)
$textRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "6308"
$ws.Cells.Item($row, 2).Value = "7/1/2025"
$ws.Cells.Item($row, 3).Value = "Guayaquil 637"
$ws.Cells.Item($row, 4).Value = "6"
$ws.Cells.Item($row, 5).Value = "807896343"
$ws.Cells.Item($row, 6).Value = "Optical Power"
$ws.Cells.Item($row, 7).Value = "Pendiente"
$ws.Cells.Item($row, 8).Value = "Picada"
$ws.Cells.Item($row, 9).Value = "1"
$ws.Cells.Item($row, 10).Value = "Cambio"
$ws.Cells.Item($row, 11).Value = "Sin equipos"
$ws.Cells.Item($row, 12).Value = "Pasante"
$ws.Cells.Item($row, 13).Value = -58.437378
$ws.Cells.Item($row, 14).Value = -34.62116

# Drop the scratch "@" formatting again so the new cells keep the sheet's
# default (unstyled) look, same as every other data row.
$textRange.Style = "Normal"
